$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Set the two new topic values for rows 32 and 33 (C column)
$ws.Range("C32").Value = "Class break due to a Program, Student's special request"
$ws.Range("C33").Value = "Interface"

# Restore the row heights of rows 36-41 to 19.5 (from 18.75)
$ws.Range("A36:A41").EntireRow.RowHeight = 19.5
